$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.843.82"
$ws.Range("E2").Value = "  -0.80%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.906.71"
$ws.Range("E3").Value = "  -0.06%  "

$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "355.53"
$ws.Range("E5").Value = "  +0.94%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "109.29"
$ws.Range("E6").Value = "  -2.40%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.562"
$ws.Range("E7").Value = "  +0.68%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.625"
$ws.Range("E9").Value = "  -0.73%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.57"
$ws.Range("E10").Value = "  -3.40%  "

$ws.Range("E11").Value = "  +1.30%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0866"
$ws.Range("E12").Value = "  +0.14%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.43"
$ws.Range("E13").Value = "  -2.27%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.68"
$ws.Range("E14").Value = "  -1.26%  "

$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.355.09"
$ws.Range("E15").Value = "  -0.44%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.889.21"
$ws.Range("E16").Value = "  -0.87%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.976"
$ws.Range("E17").Value = "  -2.88%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "51.754.73"
$ws.Range("E18").Value = "  -1.06%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.38"
$ws.Range("E19").Value = "  +1.93%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.49"
$ws.Range("E20").Value = "  -1.63%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.77"
$ws.Range("E21").Value = "  -2.14%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0973"
$ws.Range("E22").Value = "  -0.45%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.24"
$ws.Range("E23").Value = "  -1.07%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "267.37"
$ws.Range("E24").Value = "  -0.93%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.79"
$ws.Range("E25").Value = "  +0.04%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.184"
$ws.Range("E26").Value = "  +9.38%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "26.74"
$ws.Range("E27").Value = "  +0.13%  "

$ws.Range("B28").Value = "Filecoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.52"
$ws.Range("E28").Value = "  +14.30%  "

$ws.Range("B29").Value = "Dai"
$ws.Range("C29").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.12%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.105"
$ws.Range("E30").Value = "  +9.75%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "10.44"
$ws.Range("E31").Value = "  -1.55%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "37.33"
$ws.Range("E32").Value = "  -0.94%  "

$ws.Range("B33").Value = "RenderToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.13"
$ws.Range("E33").Value = "  -3.06%  "

$ws.Range("B34").Value = "Toncoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.20"
$ws.Range("E34").Value = "  -1.99%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "52.04"
$ws.Range("E35").Value = "  -2.11%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0441"
$ws.Range("E36").Value = "  -1.63%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  -0.09%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.18"
$ws.Range("E38").Value = "  -4.13%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.16"
$ws.Range("E39").Value = "  -2.53%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.98"
$ws.Range("E40").Value = "  -4.07%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.72"
$ws.Range("E41").Value = "  -3.80%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.119"
$ws.Range("E42").Value = "  +1.50%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "22.80"
$ws.Range("E43").Value = "  -4.33%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "118.90"
$ws.Range("E44").Value = "  -1.29%  "

$ws.Range("E45").Value = "  -0.88%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.47"
$ws.Range("E46").Value = "  -5.69%  "

$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.42"
$ws.Range("E47").Value = "  -3.18%  "

$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.121.11"
$ws.Range("E48").Value = "  -3.21%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.252"
$ws.Range("E49").Value = "  -4.70%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0337"
$ws.Range("E50").Value = "  +0.30%  "

$ws.Range("B51").Value = "SEI"
$ws.Range("C51").Value = "https://coinranking.com/coin/8nxCqs-uj+sei-sei"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.907"
$ws.Range("E51").Value = "  -5.77%  "
